$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "explotaciones-con-tierras-con--25-de-su-propiedad"
$ws.Range("G2").Value = "explotaciones-con-tierras-con-75-y--100-de-su-propiedad"
$ws.Range("I2").Value = "explotaciones-con-tierras-con-50-y--75-de-su-propiedad"
$ws.Range("O2").Value = "explotaciones-con-tierras-con-25-y-50-de-su-propiedad"

$ws.Range("C3").Value = "iaest-measure:explotaciones-con-tierras-con--25-de-su-propiedad"
$ws.Range("G3").Value = "iaest-measure:explotaciones-con-tierras-con-75-y--100-de-su-propiedad"
$ws.Range("I3").Value = "iaest-measure:explotaciones-con-tierras-con-50-y--75-de-su-propiedad"
$ws.Range("O3").Value = "iaest-measure:explotaciones-con-tierras-con-25-y-50-de-su-propiedad"
